$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List")

# Fill in row 8 with the new "Contains Duplicate" solution entry
$ws.Cells.Item(8, 2).Value = 217
$ws.Cells.Item(8, 3).Value = "ContainsDuplicate "
$ws.Cells.Item(8, 4).Value = "Set "
$ws.Cells.Item(8, 5).Value = "Python "

# Move the active selection to B9, as in the saved file
$ws.Range("B9").Select()
